$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "CSRD Seeds"
$ws.Range("D1").Value = "CSRD Reprasentative Terms"

# Update Average Stars values
$ws.Range("B2").Value = 3.66
$ws.Range("B3").Value = 3.81
$ws.Range("B4").Value = 3.85
$ws.Range("B5").Value = 3.71
$ws.Range("B6").Value = 3.55

$wb.Save()
